$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.895.13'
$ws.Range("E2").Value = '  +1.56%  '
$ws.Range("D3").Value = '3.253.97'
$ws.Range("E3").Value = '  +0.69%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'579.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.61%  '
$ws.Range("D6").Value = "'182.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.08%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  -1.13%  '
$ws.Range("E9").Value = '  +4.00%  '
$ws.Range("D10").Value = "'6.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.51%  '
$ws.Range("D11").Value = "'0.416"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.94%  '
$ws.Range("D12").Value = '3.816.68'
$ws.Range("E12").Value = '  +0.68%  '
$ws.Range("E13").Value = '  +0.52%  '
$ws.Range("D14").Value = "'28.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.04%  '
$ws.Range("D15").Value = '67.887.31'
$ws.Range("E15").Value = '  +1.57%  '
$ws.Range("E16").Value = '  +1.97%  '
$ws.Range("D17").Value = '3.252.29'
$ws.Range("E17").Value = '  +0.84%  '
$ws.Range("D18").Value = "'5.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.16%  '
$ws.Range("E19").Value = '  +2.62%  '
$ws.Range("D20").Value = "'379.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.70%  '
$ws.Range("D21").Value = "'7.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.08%  '
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("D23").Value = "'71.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.12%  '
$ws.Range("E24").Value = '  +2.00%  '
$ws.Range("D26").Value = "'10.05"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.45%  '
$ws.Range("E27").Value = '  +2.50%  '
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").Value = "'1.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.21%  '
$ws.Range("E30").Value = '  +2.18%  '
$ws.Range("D31").Value = "'22.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.29%  '
$ws.Range("E32").Value = '  +0.01%  '
$ws.Range("E33").Value = '  +4.37%  '
$ws.Range("E34").Value = '  +3.85%  '
$ws.Range("E35").Value = '  +5.91%  '
$ws.Range("D36").Value = "'162.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.21%  '
$ws.Range("B37").Value = 'Stacks'
$ws.Range("C37").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D37").Value = "'1.87"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.32%  '
$ws.Range("B38").Value = 'Mantle'
$ws.Range("C38").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D38").Value = "'0.844"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.94%  '
$ws.Range("D39").Value = "'26.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.53%  '
$ws.Range("E40").Value = '  +8.34%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = "'6.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.40%  '
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").Value = "'2.61"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").Value = "'25.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.01%  '
$ws.Range("D44").Value = "'346.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.14%  '
$ws.Range("D45").Value = "'41.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.85%  '
$ws.Range("E46").Value = '  +2.34%  '
$ws.Range("D47").Value = '2.631.94'
$ws.Range("E47").Value = '  -1.84%  '
$ws.Range("E48").Value = '  +3.23%  '
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("D50").Value = "'0.994"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.17%  '
$ws.Range("E51").Value = '  +2.90%  '
